# The historic-data table contained a bad/incomplete entry in row 24
# (an "ERAU" row with several blank weight fields that skewed the
# weight-estimation regression). Remove that entire row so the data
# below it shifts up, the table/used-range shrinks from R48 to R47,
# and the now-unused "ERAU" shared string drops out of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(24).Delete()

# Leave the selection where the author ended up after editing.
$ws.Range("D51").Select()
